# Nowcasts 2025Q3 update: append the 2025-08-30 revision row and refresh all figures from the latest run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$header = @("Row", "Prognose", "surveys", "production", "orders", "turnover", "financial", "labor market", "prices", "national accounts", "Revision")
for ($c = 1; $c -le $header.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $header[$c - 1]
}

# --- Row labels (A2:A12), forced to text so dates are not auto-parsed ---
$rowLabels = @("2025-03-30", "2025-04-15", "2025-04-30", "2025-05-15", "2025-05-30", "2025-06-15", "2025-06-30", "2025-07-15", "2025-07-30", "2025-08-15", "2025-08-30")
for ($i = 0; $i -lt $rowLabels.Length; $i++) {
    $cell = $ws.Cells.Item($i + 2, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $rowLabels[$i]
    $cell.Style = "Normal"
}

# --- Numeric data (B2:K12): Prognose + per-block revisions, including new row 12 ---
$data = @(
    @(0.28951790997397864, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0.2815133736964349, 0, -0.00867934578883236, -0.0005598148358634941, -0.00039016896472588184, 0.0005822541495679213, 0.0000798497642493083, 0.0008786201830608118, 0, 0.00008406921499998221),
    @(0.286862803486128, 0.0018227661922254968, 0, -0.00002347137628880709, 0.00008474784127008656, 0, 0.0002821906240403322, 0.003032259537779826, -0.00000017098273188762175, 0.00015110795339806637),
    @(0.2679146809329271, -0.00221362263038758, -0.010356518640263384, -0.0022137701167173106, -0.005971974375183251, 0.001521976083575994, -0.00020155361784646316, 0.0003692244773004711, 0, 0.00011811626632057814),
    @(0.4916688625021139, 0.2172946885511371, 0, -0.0000342419080683433, 0.0010543023650383717, 0, 0.000683433773648294, 0.004752854378015075, 0, 0.00000314440941634464),
    @(0.49821017552680824, 0, -0.006688401380593935, 0.0014134196012066435, 0.011230098141451657, 0.00023276964452411775, 0, 0.0009290143888820918, 0, -0.0005755873707762316),
    @(0.17267193400212133, -0.3228336546124008, 0, 0.00006331154439589402, -0.0035640209706762774, 0, -0.0000859787419245384, 0.0007983413386659967, 0, 0.00008375991725284582),
    @(0.040432406569817725, 0, -0.025350078806031115, -0.014096434938587869, -0.08796709070775391, -0.0025048144622745756, -0.0022216029106557843, -0.00017266691063481097, 0, 0.00007316130363446693),
    @(0.3145587803603921, 0.31109650407446776, 0, -0.0006228517202901458, 0.004917977046679238, 0, -0.00034524826166564953, -0.0007224927307180036, -0.04954093673346558, 0.009343422115566768),
    @(0.42575127339425134, 0, 0.005308790885471065, 0.017720438671829195, 0.1268206819344926, 0.0014951519468150878, -0.0001578993372511636, -0.010081711129216254, 0, -0.029912959938281303),
    @(0.3093104987787674, -0.06639871258513183, 0, 0.0011728655003356296, 0.00009604859972384022, 0, 0.00003419446692734647, -0.007224024677687497, 0, -0.04412114591965144)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value = $rowVals[$j]
    }
}

# --- Column C got a touch narrower in this run ---
# (target stored width 14.64453125; engine's ColumnWidth setter only has
# 1/6-character resolution, so 13.833333 is the closest achievable input)
$ws.Columns.Item(3).ColumnWidth = 13.833333
